# MatriceAccessi.xlsx - "Corretto diagrammi e packages"
#
# This script reproduces the content-level changes described by the
# target diff:
#   - On sheet "GestioneUtenti" the access-matrix values for the roles
#     "Gestore Prodotti" (row 4), "Gestore Utenti" (row 5) and
#     "Gestore Ordini" (row 6) are updated so that ModificaProfilo (H)
#     and VisualizzaProfilo (I) are now granted (0 -> 1).
#   - The active / selected worksheet tab moves from "GestioneOrdini"
#     back to "GestioneUtenti", and the selection on "GestioneUtenti"
#     becomes the full used range (A1:K6) instead of the single cell B3.

$wb = $excel.ActiveWorkbook

$wsUtenti  = $wb.Worksheets.Item("GestioneUtenti")
$wsOrdini  = $wb.Worksheets.Item("GestioneOrdini")

# --- Update the access matrix values -------------------------------------
# Row 4 = Gestore Prodotti, Row 5 = Gestore Utenti, Row 6 = Gestore Ordini
# Column H = ModificaProfilo, Column I = VisualizzaProfilo
$wsUtenti.Range("H4").Value = 1
$wsUtenti.Range("I4").Value = 1

$wsUtenti.Range("H5").Value = 1
$wsUtenti.Range("I5").Value = 1

$wsUtenti.Range("H6").Value = 1
$wsUtenti.Range("I6").Value = 1

# --- Switch the active tab back to "GestioneUtenti" -----------------------
# Selecting the range first on the (currently) active sheet makes sure its
# own selection/highlight changes before we switch the active tab.
$wsOrdini.Activate() | Out-Null
$wsOrdini.Range("F6").Select() | Out-Null

$wsUtenti.Activate() | Out-Null
$wsUtenti.Range("A1:K6").Select() | Out-Null
